$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Logic" operator characters to the secondary (Trig-overlay)
#     keypad columns, mirroring the existing "(" / ")" slots that used to sit
#     there. L4/M4 belong to the first duplicated keypad block, W4/X4 to the
#     second; L5 becomes the new "Logic" label (previously "Trig").
$ws.Range("L4").Value = "<"
$ws.Range("M4").Value = ">"
$ws.Range("W4").Value = [char]0x2264   # "<=" (less-than-or-equal)
$ws.Range("X4").Value = [char]0x2265   # ">=" (greater-than-or-equal)
$ws.Range("L5").Value = "Logic"

# --- Re-apply the highlighted keypad-button look (font/fill/border) that the
#     rest of row 4's lettered keys already use, so the new cells match the
#     surrounding "special key" styling instead of the plain look they
#     inherited from the old "(" / ")" cells.
$ws.Range("AF1").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("W4").PasteSpecial(-4122)
$ws.Range("X4").PasteSpecial(-4122)

# --- L5 ("Logic") picks up the same style family as its row-5 neighbor
#     (the already-highlighted "FUNC\nRCLL"-style header cell).
$ws.Range("AA2").Copy()
$ws.Range("L5").PasteSpecial(-4122)

# --- Move the saved selection/viewport: scroll back to column A and select
#     R12 instead of the old M9 selection left over from editing.
$ws.Range("R12").Select()
